$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Staff" to "data"
$ws.Name = "data"

# Update the "Cepas Can Number" value for the last row (F4)
$ws.Range("F4").Value = "Poekoas"

# Update the "Printed Date" value for the last row (G4) - keep it as text
$ws.Range("G4").Value = "08/22/2023 04:15:53 PM"
